$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the question text in E3: add ":pflicht" qualifier to both parts
$ws.Range("E3").Value = "Ticketart:dropdown(Bus,Zug,U-Bahn):pflicht;Häufigkeit:dropdown(Täglich,Wöchentlich,Selten):pflicht"

# Update the active selection to E3
$ws.Range("E3").Select()
